$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "contingencies with rene fine": two new line contingency rows (line7, line8)
# are inserted into the table, pushing the existing extr1..extr8 rows down by
# two positions, and two brand-new extr7/extr8 rows are appended at the
# bottom (rows 16 and 17). The underlying from_bus/to_bus numbers for each
# extr* label are unchanged from before - only their row position and some
# in_service flags change - while line7/line8 carry entirely new data.

# Final table (rows 2-17), column order: A=id, B=name, C=from_bus, D=to_bus, E=in_service
$rows = @(
    @(2,  0, "line1", 7,  9,  $true),
    @(3,  1, "line2", 9,  8,  $false),
    @(4,  2, "line3", 8,  10, $true),
    @(5,  3, "line4", 8,  11, $true),
    @(6,  4, "line5", 10, 5,  $true),
    @(7,  5, "line6", 12, 8,  $true),
    @(8,  6, "line7", 14, 11, $true),
    @(9,  7, "line8", 16, 9,  $true),
    @(10, 8, "extr1", 5,  12, $true),
    @(11, 9, "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $false),
    @(13, 11, "extr4", 7,  8,  $true),
    @(14, 12, "extr5", 9,  11, $true),
    @(15, 13, "extr6", 7,  11, $true),
    @(16, 14, "extr7", 5,  7,  $false),
    @(17, 15, "extr8", 8,  5,  $false)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# Columns A for the two brand new rows (16, 17) need the same formatting
# (bold, centered, bordered) as the rest of the A column - copy it over from
# an existing formatted cell so we reuse the workbook's existing style
# instead of minting new ones.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
